$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.049.46'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.358.96'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.679'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.75'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.40'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.20%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.593'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +9.86%  '

$ws.Range("E10").Value = '  +1.33%  '

$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.06'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +13.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.27'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +9.33%  '

$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.705.67'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.64'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.903'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.349.27'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.927.70'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '76.99'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.28'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +24.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.68'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.78%  '

$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.71'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("E29").Value = '  -1.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.75'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.33'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.56%  '

$ws.Range("E32").Value = '  -2.29%  '

$ws.Range("E33").Value = '  +3.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0762'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.24'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.43'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.75'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.82%  '

$ws.Range("E38").Value = '  -2.81%  '

$ws.Range("E39").Value = '  -1.49%  '

$ws.Range("E40").Value = '  +4.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.115'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +16.98%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.15'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.75%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.21'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.204'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +12.68%  '

$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.72'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.91'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +10.26%  '

$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.24'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.48'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +7.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.06%  '

$ws.Range("E51").Value = '  +2.67%  '
